$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (columns A-D) ----
$ws.Cells.Item(1, 1).Value = "CÓDIGO OC"
$ws.Cells.Item(1, 2).Value = "REFERENCIA"
$ws.Cells.Item(1, 3).Value = "CANTIDAD DESPACHADA"
$ws.Cells.Item(1, 4).Value = "FECHA DESPACHO"

# ---- Data row 2 ----
$ws.Cells.Item(2, 1).Value = "OCAM-63-2021-115"
$ws.Cells.Item(2, 2).Value = "CI1-U"
$ws.Cells.Item(2, 3).Value = 139
$ws.Cells.Item(2, 4).Value = "07/ene/2022"

# ---- Data row 4 (entered before row 3) ----
$ws.Cells.Item(4, 1).Value = "OCAM-57-2021-115"
$ws.Cells.Item(4, 2).Value = "CIH1-L"
$ws.Cells.Item(4, 3).Value = 51
$ws.Cells.Item(4, 4).Value = "06/ene/2021"

# ---- Data row 3 ----
$ws.Cells.Item(3, 1).Value = "OCAM-57-2021-115"
$ws.Cells.Item(3, 2).Value = "CIH1-XXL"
$ws.Cells.Item(3, 3).Value = 32
$ws.Cells.Item(3, 4).Value = "11/ene/2022"

# ---- Remisión column (E), added last ----
$ws.Cells.Item(1, 5).Value = "# REMISIÓN"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(4, 5).Value = 2

# ---- Header style: bold white text, blue fill, centered ----
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 12611584
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# ---- Data style: centered ----
$dataRange = $ws.Range("A2:E4")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 15.0
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668
$ws.Columns.Item(5).ColumnWidth = 11.5

# ---- Selection ----
$ws.Range("E2").Select()
